$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All edited cells hold text (coin prices/links/percent deltas are
# displayed as formatted strings, not live numbers), so force the
# "Price" column to a text format before writing -- otherwise plain
# decimal-looking values like "601.95" would be auto-coerced to numbers.
$priceCells = @(
    'D2', 'D3', 'D5', 'D6', 'D8', 'D9', 'D11', 'D12', 'D13', 'D14', 'D15', 'D16', 'D17', 'D18', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D27', 'D28', 'D29', 'D31', 'D32', 'D35', 'D36', 'D37', 'D38', 'D39', 'D42', 'D43', 'D44', 'D45', 'D46', 'D48', 'D51'
)
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '66.807.98'
$ws.Range('E2').Value = '  +3.19%  '
$ws.Range('D3').Value = '3.210.66'
$ws.Range('E3').Value = '  +2.36%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '601.95'
$ws.Range('E5').Value = '  +3.79%  '
$ws.Range('D6').Value = '158.22'
$ws.Range('E6').Value = '  +7.64%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').Value = '3.210.72'
$ws.Range('D9').Value = '0.554'
$ws.Range('E9').Value = '  +5.72%  '
$ws.Range('E10').Value = '  +2.35%  '
$ws.Range('D11').Value = '6.00'
$ws.Range('E11').Value = '  -2.25%  '
$ws.Range('D12').Value = '0.516'
$ws.Range('E12').Value = '  +3.88%  '
$ws.Range('D13').Value = '0.0000268'
$ws.Range('E13').Value = '  +1.59%  '
$ws.Range('D14').Value = '39.33'
$ws.Range('E14').Value = '  +6.34%  '
$ws.Range('D15').Value = '3.738.21'
$ws.Range('E15').Value = '  +2.37%  '
$ws.Range('D16').Value = '66.835.31'
$ws.Range('E16').Value = '  +3.15%  '
$ws.Range('D17').Value = '7.47'
$ws.Range('E17').Value = '  +4.99%  '
$ws.Range('D18').Value = '3.213.26'
$ws.Range('E18').Value = '  +2.34%  '
$ws.Range('E19').Value = '  +0.99%  '
$ws.Range('D20').Value = '517.86'
$ws.Range('E20').Value = '  +3.94%  '
$ws.Range('D21').Value = '15.43'
$ws.Range('E21').Value = '  +1.20%  '
$ws.Range('D22').Value = '0.744'
$ws.Range('E22').Value = '  +4.73%  '
$ws.Range('D23').Value = '8.21'
$ws.Range('E23').Value = '  +6.30%  '
$ws.Range('D24').Value = '15.06'
$ws.Range('E24').Value = '  +0.51%  '
$ws.Range('D25').Value = '85.40'
$ws.Range('E25').Value = '  +1.49%  '
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').Value = '9.38'
$ws.Range('E27').Value = '  +3.54%  '
$ws.Range('D28').Value = '3.03'
$ws.Range('E28').Value = '  +4.38%  '
$ws.Range('D29').Value = '2.44'
$ws.Range('E29').Value = '  +12.23%  '
$ws.Range('E30').Value = '  +10.51%  '
$ws.Range('D31').Value = '7.06'
$ws.Range('E31').Value = '  +10.53%  '
$ws.Range('D32').Value = '28.34'
$ws.Range('E32').Value = '  +3.44%  '
$ws.Range('E33').Value = '  +2.07%  '
$ws.Range('E34').Value = '  +0.10%  '
$ws.Range('D35').Value = '6.61'
$ws.Range('E35').Value = '  +2.71%  '
$ws.Range('D36').Value = '527.00'
$ws.Range('E36').Value = '  +12.71%  '
$ws.Range('D37').Value = '55.07'
$ws.Range('E37').Value = '  +0.56%  '
$ws.Range('D38').Value = '0.0906'
$ws.Range('E38').Value = '  +1.87%  '
$ws.Range('D39').Value = '0.0426'
$ws.Range('E39').Value = '  +2.53%  '
$ws.Range('E40').Value = '  +9.89%  '
$ws.Range('E41').Value = '  +2.15%  '
$ws.Range('D42').Value = '8.95'
$ws.Range('E42').Value = '  +2.81%  '
$ws.Range('B43').Value = 'PEPE'
$ws.Range('C43').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D43').Value = '0.0₃0697'
$ws.Range('E43').Value = '  +16.99%  '
$ws.Range('B44').Value = 'TheGraph'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D44').Value = '0.305'
$ws.Range('E44').Value = '  +8.82%  '
$ws.Range('D45').Value = '2.50'
$ws.Range('E45').Value = '  +3.29%  '
$ws.Range('D46').Value = '2.896.05'
$ws.Range('E46').Value = '  -2.45%  '
$ws.Range('E47').Value = '  +2.82%  '
$ws.Range('D48').Value = '2.41'
$ws.Range('E48').Value = '  +8.52%  '
$ws.Range('E49').Value = '  +3.60%  '
$ws.Range('D51').Value = '2.67'
$ws.Range('E51').Value = '  +10.43%  '
